$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set value "x" in column J for these rows (already has "x" in G,H,I)
$rowsJ = @(6,7,8,10,11,12,13,15)
foreach ($r in $rowsJ) {
    $ws.Range("J$r").Value = "x"
}

# Row 17 and 18: I and J columns both get "x"
$ws.Range("I17").Value = "x"
$ws.Range("J17").Value = "x"
$ws.Range("I18").Value = "x"
$ws.Range("J18").Value = "x"

# Update the active selection to J18
$ws.Range("J18").Select()
